$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.455.46"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.836.75"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +1.10%  "
$ws.Range("D5").Value = "'314.41"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "'0.4738"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("D8").Value = "'0.3692"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "'0.07460"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("D10").Value = "'0.8854"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "1.890.14"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").Value = "'0.07341"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").Value = "'5.448"
$ws.Range("D15").Value = "'93.23"
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "'6.580"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "'0.000008819"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "'14.81"
$ws.Range("E20").Value = "  +1.08%  "
$ws.Range("D21").Value = "27.466.78"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").Value = "'5.319"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").Value = "'10.69"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "2.105.75"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").Value = "'1.902"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'152.12"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "'18.65"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D28").Value = "'2.143"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'5.256"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'117.69"
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("D31").Value = "'0.08999"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").Value = "'0.7571"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "'1.179"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").Value = "'4.551"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").Value = "'2.945"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("D38").Value = "'0.05343"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").Value = "'7.319"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "'2.403"
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("D43").Value = "'0.5332"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'0.1661"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "'8.515"
$ws.Range("E45").Value = "  +1.04%  "
$ws.Range("D46").Value = "'0.4908"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "'10.55"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "'104.95"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").Value = "'1.677"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").Value = "'0.06307"
